$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '61.885.07'
$ws.Range('E2').Value = '  +4.64%  '
$ws.Range('D3').Value = '3.080.61'
$ws.Range('E3').Value = '  +3.08%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '579.64'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.06%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.22'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.79%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.069.43'
$ws.Range('E8').Value = '  +3.15%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.527'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.30%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.140'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.37%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.51'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +8.26%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.468'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.64%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000240'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +4.66%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '35.23'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +4.69%  '
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').Value = '3.594.85'
$ws.Range('E16').Value = '  +3.20%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.26'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.12%  '
$ws.Range('D18').Value = '3.082.12'
$ws.Range('E18').Value = '  +3.00%  '
$ws.Range('D19').Value = '61.849.12'
$ws.Range('E19').Value = '  +4.63%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '448.37'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.89%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.91'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.77%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.729'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.04%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.46'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +5.42%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.81'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.71%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '82.10'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.94%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.25'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +5.66%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  +4.97%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.08'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +4.98%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.71'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +9.71%  '
$ws.Range('E32').Value = '  +14.72%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '26.73'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +4.27%  '
$ws.Range('E34').Value = '  +4.40%  '
$ws.Range('D35').Value = '0.0₃0797'
$ws.Range('E35').Value = '  +3.36%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.04'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +3.93%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.18'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +5.66%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '50.37'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.72%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.99'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +9.59%  '
$ws.Range('E40').Value = '  +2.17%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '429.58'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +7.86%  '
$ws.Range('E42').Value = '  +6.00%  '
$ws.Range('D43').Value = '2.794.26'
$ws.Range('E43').Value = '  +1.72%  '
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.269'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +6.99%  '
$ws.Range('B46').Value = 'Arweave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '35.33'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +9.13%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.09'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +4.97%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '124.87'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('E50').Value = '  +1.37%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '24.00'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.33%  '
